# Update "想去人数" (want-to-go count, column F) figures across the four
# sheets of the 广州-漫展信息 workbook, per the latest scrape refresh
# (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 84
$ws.Range("F5").Value  = 168
$ws.Range("F6").Value  = 395
$ws.Range("F7").Value  = 178
$ws.Range("F9").Value  = 1031
$ws.Range("F10").Value = 340
$ws.Range("F11").Value = 177
$ws.Range("F12").Value = 43
$ws.Range("F14").Value = 363
$ws.Range("F15").Value = 349
$ws.Range("F16").Value = 761
$ws.Range("F19").Value = 256
$ws.Range("F20").Value = 67
$ws.Range("F21").Value = 975
$ws.Range("F22").Value = 434
$ws.Range("F23").Value = 246
$ws.Range("F25").Value = 360
$ws.Range("F28").Value = 456

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 356
$ws.Range("F6").Value  = 40
$ws.Range("F7").Value  = 278
$ws.Range("F11").Value = 145

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 343

# Sheet "全部类型" (All types) - union of the other sheets
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 343
$ws.Range("F6").Value  = 84
$ws.Range("F7").Value  = 168
$ws.Range("F8").Value  = 395
$ws.Range("F9").Value  = 178
$ws.Range("F11").Value = 1031
$ws.Range("F12").Value = 340
$ws.Range("F13").Value = 177
$ws.Range("F15").Value = 43
$ws.Range("F16").Value = 356
$ws.Range("F19").Value = 363
$ws.Range("F20").Value = 40
$ws.Range("F21").Value = 278
$ws.Range("F22").Value = 349
$ws.Range("F23").Value = 761
$ws.Range("F26").Value = 256
$ws.Range("F27").Value = 67
$ws.Range("F28").Value = 975
$ws.Range("F29").Value = 434
$ws.Range("F32").Value = 246
$ws.Range("F34").Value = 360
$ws.Range("F36").Value = 145
$ws.Range("F40").Value = 456
